$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column (C) for rows 2-56 is updated
# from serial date 45175 (2023-09-06) to 45183 (2023-09-14).
$ws.Range("C2:C56").Value = 45183
